# --- NYPD CompStat weekly refresh: new crime data collected ---
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Report header: volume/number + week-covering date range (rich-text runs
# in the shared string all share identical formatting, so plain text
# assignment reproduces the same rendered result).
$ws.Range("A8").Value = "Volume 31   Number  33"
$ws.Range("C9").Value = "Report Covering the Week  8/12/2024  Through  8/18/2024"

# Crime-complaints table (rows 14-31): refreshed weekly/28-day/YTD counts
# and their derived percent-change figures.
$ws.Range("D14").Value = 2
$ws.Range("G14").Value = 3
$ws.Range("H14").Value = -66.666666666666
$ws.Range("J14").Value = 14
$ws.Range("K14").Value = 7.142857142857
$ws.Range("N14").Value = -72.727272727272
$ws.Range("C15").Value = 3
$ws.Range("D15").Value = 2
$ws.Range("E15").Value = 50
$ws.Range("G15").Value = 10
$ws.Range("H15").Value = 20
$ws.Range("I15").Value = 95
$ws.Range("J15").Value = 76
$ws.Range("K15").Value = 25
$ws.Range("L15").Value = -15.929203539823
$ws.Range("M15").Value = 55.737704918032
$ws.Range("N15").Value = -33.098591549295
$ws.Range("C16").Value = 39
$ws.Range("D16").Value = 40
$ws.Range("E16").Value = -2.5
$ws.Range("F16").Value = 161
$ws.Range("G16").Value = 151
$ws.Range("H16").Value = 6.622516556291
$ws.Range("I16").Value = 1029
$ws.Range("J16").Value = 1138
$ws.Range("K16").Value = -9.57820738137
$ws.Range("L16").Value = -21.689497716895
$ws.Range("M16").Value = 32.603092783505
$ws.Range("N16").Value = -84.41145281018
$ws.Range("C17").Value = 37
$ws.Range("D17").Value = 52
$ws.Range("E17").Value = -28.846153846153
$ws.Range("F17").Value = 187
$ws.Range("G17").Value = 197
$ws.Range("H17").Value = -5.076142131979
$ws.Range("I17").Value = 1466
$ws.Range("J17").Value = 1412
$ws.Range("K17").Value = 3.824362606232
$ws.Range("L17").Value = 13.030069390902
$ws.Range("M17").Value = 82.565379825653
$ws.Range("N17").Value = -32.967535436671
$ws.Range("C18").Value = 27
$ws.Range("D18").Value = 51
$ws.Range("E18").Value = -47.058823529411
$ws.Range("F18").Value = 109
$ws.Range("G18").Value = 168
$ws.Range("H18").Value = -35.119047619047
$ws.Range("I18").Value = 1162
$ws.Range("J18").Value = 1348
$ws.Range("K18").Value = -13.798219584569
$ws.Range("L18").Value = -38.906414300736
$ws.Range("M18").Value = 2.288732394366
$ws.Range("N18").Value = -84.398496240601
$ws.Range("C19").Value = 231
$ws.Range("D19").Value = 224
$ws.Range("E19").Value = 3.125
$ws.Range("F19").Value = 846
$ws.Range("G19").Value = 910
$ws.Range("H19").Value = -7.032967032967
$ws.Range("I19").Value = 6450
$ws.Range("J19").Value = 7262
$ws.Range("K19").Value = -11.181492701735
$ws.Range("L19").Value = -9.891031014249
$ws.Range("M19").Value = -0.416859657248
$ws.Range("N19").Value = -68.496629872032
$ws.Range("C20").Value = 14
$ws.Range("D20").Value = 17
$ws.Range("E20").Value = -17.647058823529
$ws.Range("F20").Value = 39
$ws.Range("H20").Value = -36.065573770491
$ws.Range("I20").Value = 278
$ws.Range("J20").Value = 381
$ws.Range("K20").Value = -27.034120734908
$ws.Range("L20").Value = -37.807606263982
$ws.Range("M20").Value = 15.352697095435
$ws.Range("N20").Value = -92.661034846884
$ws.Range("C21").Value = 351
$ws.Range("D21").Value = 388
$ws.Range("E21").Value = -9.536082474226
$ws.Range("F21").Value = 1355
$ws.Range("G21").Value = 1500
$ws.Range("H21").Value = -9.666666666666
$ws.Range("I21").Value = 10495
$ws.Range("J21").Value = 11631
$ws.Range("K21").Value = -9.767001977473
$ws.Range("L21").Value = -14.340515834149
$ws.Range("M21").Value = 10.450431488107
$ws.Range("N21").Value = -74.210590981693
$ws.Range("C22").Value = 14
$ws.Range("D22").Value = 15
$ws.Range("E22").Value = -6.666666666666
$ws.Range("F22").Value = 49
$ws.Range("G22").Value = 49
$ws.Range("H22").Value = 0
$ws.Range("I22").Value = 389
$ws.Range("J22").Value = 420
$ws.Range("K22").Value = -7.380952380952
$ws.Range("L22").Value = -8.037825059101
$ws.Range("M22").Value = 18.237082066869
$ws.Range("C23").Value = 12
$ws.Range("D23").Value = 6
$ws.Range("E23").Value = 100
$ws.Range("G23").Value = 34
$ws.Range("H23").Value = -20.588235294117
$ws.Range("I23").Value = 247
$ws.Range("J23").Value = 259
$ws.Range("K23").Value = -4.633204633204
$ws.Range("L23").Value = -15.120274914089
$ws.Range("M23").Value = 8.810572687224
$ws.Range("C24").Value = 541
$ws.Range("D24").Value = 476
$ws.Range("E24").Value = 13.655462184873
$ws.Range("F24").Value = 1875
$ws.Range("G24").Value = 1771
$ws.Range("H24").Value = 5.872388481084
$ws.Range("I24").Value = 13870
$ws.Range("J24").Value = 12979
$ws.Range("K24").Value = 6.864935665305
$ws.Range("L24").Value = 1.634058767494
$ws.Range("M24").Value = 31.519059358998
$ws.Range("C25").Value = 461
$ws.Range("D25").Value = 411
$ws.Range("E25").Value = 12.165450121654
$ws.Range("F25").Value = 1584
$ws.Range("G25").Value = 1457
$ws.Range("H25").Value = 8.716540837336
$ws.Range("I25").Value = 11763
$ws.Range("J25").Value = 10810
$ws.Range("K25").Value = 8.815911193339
$ws.Range("L25").Value = 1.065383624022
$ws.Range("C26").Value = 108
$ws.Range("D26").Value = 114
$ws.Range("E26").Value = -5.263157894736
$ws.Range("F26").Value = 415
$ws.Range("G26").Value = 409
$ws.Range("H26").Value = 1.466992665036
$ws.Range("I26").Value = 3144
$ws.Range("J26").Value = 3150
$ws.Range("K26").Value = -0.190476190476
$ws.Range("L26").Value = 7.230559345156
$ws.Range("M26").Value = 37.232649498035
$ws.Range("C27").Value = 3
$ws.Range("D27").Value = 3
$ws.Range("E27").Value = 0
$ws.Range("F27").Value = 17
$ws.Range("G27").Value = 18
$ws.Range("H27").Value = -5.555555555555
$ws.Range("I27").Value = 148
$ws.Range("J27").Value = 136
$ws.Range("K27").Value = 8.823529411764
$ws.Range("L27").Value = -16.38418079096
$ws.Range("C28").Value = 16
$ws.Range("D28").Value = 24
$ws.Range("E28").Value = -33.333333333333
$ws.Range("F28").Value = 93
$ws.Range("G28").Value = 76
$ws.Range("H28").Value = 22.368421052631
$ws.Range("I28").Value = 594
$ws.Range("J28").Value = 576
$ws.Range("K28").Value = 3.125
$ws.Range("L28").Value = -1.655629139072
$ws.Range("F29").Value = 4
$ws.Range("G29").Value = 3
$ws.Range("H29").Value = 33.333333333333
$ws.Range("I29").Value = 25
$ws.Range("K29").Value = -10.714285714285
$ws.Range("L29").Value = -32.432432432432
$ws.Range("M29").Value = 4.166666666666
$ws.Range("N29").Value = -73.118279569892
$ws.Range("F30").Value = 3
$ws.Range("G30").Value = 3
$ws.Range("H30").Value = 0
$ws.Range("I30").Value = 21
$ws.Range("K30").Value = -8.695652173913
$ws.Range("L30").Value = -34.375
$ws.Range("M30").Value = 16.666666666666
$ws.Range("N30").Value = -74.074074074074
$ws.Range("D31").Value = 3
$ws.Range("E31").Value = -100
$ws.Range("G31").Value = 8
$ws.Range("H31").Value = 0
$ws.Range("I31").Value = 95
$ws.Range("J31").Value = 79
$ws.Range("K31").Value = 20.253164556962
$ws.Range("L31").Value = -15.929203539823

# Murder count under "Hate Crimes" (C31) drops to the table's textual
# "0" placeholder (same convention already used at C14), not the number 0.
# Copying the already-typed source cell (instead of assigning .Value)
# keeps the exact text-shared-string + cell style pairing intact.
$ws.Range("C14").Copy($ws.Range("C31"))
